# Regen save_data to use K instead of Strike#: update the K (strikeouts)
# column values (column G) for each game row in the pitching log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new K value (column G), taken from the regenerated save data.
$kValues = [ordered]@{
    2  = 2
    3  = 2
    4  = 1
    5  = 0
    6  = 0
    7  = 1
    8  = 2
    9  = 1
    10 = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 2
    22 = 4
    23 = 1
    24 = 2
    25 = 2
    26 = 0
    28 = 0
    29 = 0
    30 = 0
    31 = 2
    32 = 2
    33 = 1
    34 = 0
    35 = 0
    36 = 1
    37 = 1
    38 = 0
    39 = 2
    40 = 2
    41 = 1
    42 = 5
    43 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
